$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matching the inlineStr type in the source file).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "33.755.81"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.765.71"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "224.13"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").Value = "0.545"
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("D8").Value = "31.98"
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "0.0686"
$ws.Range("E10").Value = "  -3.13%  "
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "2.020.82"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "11.27"
$ws.Range("E13").Value = "  +7.25%  "
$ws.Range("D14").Value = "1.770.67"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "33.728.51"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").Value = "66.51"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "237.88"
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "10.55"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").Value = "159.29"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "16.09"
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "7.01"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "1.22"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("D32").Value = "3.59"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("D33").Value = "3.49"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("D35").Value = "1.376.70"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "0.651"
$ws.Range("E36").Value = "  +2.35%  "
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").Value = "0.0184"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").Value = "2.23"
$ws.Range("E39").Value = "  +5.82%  "
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "13.65"
$ws.Range("E41").Value = "  +16.37%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.903"
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "77.38"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.65"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("E45").Value = "  +4.38%  "
$ws.Range("E46").Value = "  +14.83%  "
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").Value = "107.35"
$ws.Range("E48").Value = "  +3.19%  "
$ws.Range("D49").Value = "5.80"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").Value = "1.920.05"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").Value = "  +0.57%  "

# Restore default (Normal) style on cells where we temporarily forced a text number format,
# so no stray style attribute is left behind on the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
